$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (column BO = 67), copy formatting from BN1
$ws.Range("BO1").Value = "12-sep"
$ws.Range("BO1").NumberFormat = $ws.Range("BN1").NumberFormat

# New data values for each row, copying formatting from the BN column
$values = 15, 14, 12, 13, 11, 16, 22, 10, 11, 12

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $srcCell = $ws.Range("BN$row")
    $dstCell = $ws.Range("BO$row")
    $dstCell.Value = $values[$i]
    $dstCell.HorizontalAlignment = $srcCell.HorizontalAlignment
    $dstCell.NumberFormat = $srcCell.NumberFormat
}

# Update the selection to match the saved selection in the workbook
$ws.Range("BR5").Select()
